$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): columns B..BC ---
$header = New-Object 'object[,]' 1,54
$header[0,0] = "country_index"
$header[0,1] = "Hult_Team_Regions"
$header[0,2] = "country_name"
$header[0,3] = "country_code"
$header[0,4] = "income_group"
$header[0,5] = "access_to_electricity_pop"
$header[0,6] = "access_to_electricity_rural"
$header[0,7] = "access_to_electricity_urban"
$header[0,8] = "CO2_emissions_per_capita"
$header[0,9] = "compulsory_edu_yrs"
$header[0,10] = "pct_female_employment"
$header[0,11] = "pct_male_employment"
$header[0,12] = "pct_agriculture_employment"
$header[0,13] = "gni_index"
$header[0,14] = "pct_industry_employment"
$header[0,15] = "pct_services_employment"
$header[0,16] = "exports_pct_gdp"
$header[0,17] = "fdi_pct_gdp"
$header[0,18] = "gdp_usd"
$header[0,19] = "gdp_growth_pct"
$header[0,20] = "incidence_hiv"
$header[0,21] = "internet_usage_pct"
$header[0,22] = "child_mortality_per_1k"
$header[0,23] = "avg_air_pollution"
$header[0,24] = "women_in_parliament"
$header[0,25] = "unemployment_pct"
$header[0,26] = "urban_population_pct"
$header[0,27] = "urban_population_growth_pct"
$header[0,28] = "m_income_group"
$header[0,29] = "m_access_to_electricity_pop"
$header[0,30] = "m_access_to_electricity_rural"
$header[0,31] = "m_access_to_electricity_urban"
$header[0,32] = "m_CO2_emissions_per_capita"
$header[0,33] = "m_compulsory_edu_yrs"
$header[0,34] = "m_pct_female_employment"
$header[0,35] = "m_pct_male_employment"
$header[0,36] = "m_pct_agriculture_employment"
$header[0,37] = "m_pct_industry_employment"
$header[0,38] = "m_pct_services_employment"
$header[0,39] = "m_exports_pct_gdp"
$header[0,40] = "m_fdi_pct_gdp"
$header[0,41] = "m_gdp_usd"
$header[0,42] = "m_gdp_growth_pct"
$header[0,43] = "m_incidence_hiv"
$header[0,44] = "m_internet_usage_pct"
$header[0,45] = "m_homicides_per_100k"
$header[0,46] = "m_adult_literacy_pct"
$header[0,47] = "m_child_mortality_per_1k"
$header[0,48] = "m_avg_air_pollution"
$header[0,49] = "m_women_in_parliament"
$header[0,50] = "m_tax_revenue_pct_gdp"
$header[0,51] = "m_unemployment_pct"
$header[0,52] = "m_urban_population_pct"
$header[0,53] = "m_urban_population_growth_pct"
$ws.Range("B1:BC1").Value = $header

# Apply the existing header style (s="1", bold+border+center) to the newly added header cells AL1:BC1
$ws.Range("AK1").Copy()
$ws.Range("AL1:BC1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2 ---
$row2 = New-Object 'object[,]' 1,54
$row2[0,0] = 32
$row2[0,1] = "Central Africa 1"
$row2[0,2] = "Burundi"
$row2[0,3] = "BDI"
$row2[0,4] = "Low income"
$row2[0,5] = 7
$row2[0,6] = 2
$row2[0,7] = 52.1
$row2[0,8] = 0.044485376
$row2[0,9] = 9.751295336787564
$row2[0,10] = 12.70699978
$row2[0,11] = 3.730000019
$row2[0,12] = 91.10299683
$row2[0,13] = 290
$row2[0,14] = 2.549000025
$row2[0,15] = 6.34800005
$row2[0,16] = 7.76818252
$row2[0,17] = 2.642421428
$row2[0,18] = 3093647227
$row2[0,19] = 4.660918184
$row2[0,20] = 0.04
$row2[0,21] = 1.38
$row2[0,22] = 77.90000000000001
$row2[0,23] = 47.08506861
$row2[0,24] = 30.5
$row2[0,25] = 1.570000052
$row2[0,26] = 11.761
$row2[0,27] = 5.480228006
$row2[0,28] = 0
$row2[0,29] = 0
$row2[0,30] = 0
$row2[0,31] = 0
$row2[0,32] = 0
$row2[0,33] = 1
$row2[0,34] = 0
$row2[0,35] = 0
$row2[0,36] = 0
$row2[0,37] = 0
$row2[0,38] = 0
$row2[0,39] = 0
$row2[0,40] = 0
$row2[0,41] = 0
$row2[0,42] = 0
$row2[0,43] = 0
$row2[0,44] = 0
$row2[0,45] = 0
$row2[0,46] = 0
$row2[0,47] = 0
$row2[0,48] = 0
$row2[0,49] = 0
$row2[0,50] = 1
$row2[0,51] = 0
$row2[0,52] = 0
$row2[0,53] = 0
$ws.Range("B2:BC2").Value = $row2

# --- Row 3 ---
$row3 = New-Object 'object[,]' 1,54
$row3[0,0] = 33
$row3[0,1] = "Central Africa 1"
$row3[0,2] = "Cabo Verde"
$row3[0,3] = "CPV"
$row3[0,4] = "Lower middle income"
$row3[0,5] = 87.87672424
$row3[0,6] = 81.25119979
$row3[0,7] = 91.46945952999999
$row3[0,8] = 0.9334032370000001
$row3[0,9] = 10
$row3[0,10] = 17.28199959
$row3[0,11] = 4.796000004
$row3[0,12] = 68.65499878
$row3[0,13] = 3310
$row3[0,14] = 6.922999859
$row3[0,15] = 24.42200089
$row3[0,16] = 40.36113967
$row3[0,17] = 9.718943229000001
$row3[0,18] = 1858121723
$row3[0,19] = 0.611212666
$row3[0,20] = 0.03
$row3[0,21] = 40.26
$row3[0,22] = 23
$row3[0,23] = 42.09588617
$row3[0,24] = 20.8
$row3[0,25] = 10.34799957
$row3[0,26] = 64.84
$row3[0,27] = 2.315808738
$row3[0,28] = 0
$row3[0,29] = 0
$row3[0,30] = 0
$row3[0,31] = 0
$row3[0,32] = 0
$row3[0,33] = 0
$row3[0,34] = 0
$row3[0,35] = 0
$row3[0,36] = 0
$row3[0,37] = 0
$row3[0,38] = 0
$row3[0,39] = 0
$row3[0,40] = 0
$row3[0,41] = 0
$row3[0,42] = 0
$row3[0,43] = 1
$row3[0,44] = 0
$row3[0,45] = 0
$row3[0,46] = 1
$row3[0,47] = 0
$row3[0,48] = 0
$row3[0,49] = 0
$row3[0,50] = 1
$row3[0,51] = 0
$row3[0,52] = 0
$row3[0,53] = 0
$ws.Range("B3:BC3").Value = $row3

# --- Row 4 ---
$row4 = New-Object 'object[,]' 1,54
$row4[0,0] = 44
$row4[0,1] = "Central Africa 1"
$row4[0,2] = "Comoros"
$row4[0,3] = "COM"
$row4[0,4] = "Low income"
$row4[0,5] = 72.91772460999999
$row4[0,6] = 66.79169892
$row4[0,7] = 88.52062988
$row4[0,8] = 0.202814119
$row4[0,9] = 6
$row4[0,10] = 18.89999962
$row4[0,11] = 9.661999701999999
$row4[0,12] = 54.98300171
$row4[0,13] = 830
$row4[0,14] = 15.57699966
$row4[0,15] = 29.44099998
$row4[0,16] = 16.87471351
$row4[0,17] = 0.72227588
$row4[0,18] = 647720707.1
$row4[0,19] = 2.061639469
$row4[0,20] = 0.01
$row4[0,21] = 6.98
$row4[0,22] = 78.3
$row4[0,23] = 17.11820358
$row4[0,24] = 3
$row4[0,25] = 4.361999989
$row4[0,26] = 28.193
$row4[0,27] = 2.699437729
$row4[0,28] = 0
$row4[0,29] = 0
$row4[0,30] = 0
$row4[0,31] = 0
$row4[0,32] = 0
$row4[0,33] = 0
$row4[0,34] = 0
$row4[0,35] = 0
$row4[0,36] = 0
$row4[0,37] = 0
$row4[0,38] = 0
$row4[0,39] = 0
$row4[0,40] = 0
$row4[0,41] = 0
$row4[0,42] = 0
$row4[0,43] = 0
$row4[0,44] = 0
$row4[0,45] = 1
$row4[0,46] = 1
$row4[0,47] = 0
$row4[0,48] = 0
$row4[0,49] = 0
$row4[0,50] = 1
$row4[0,51] = 0
$row4[0,52] = 0
$row4[0,53] = 0
$ws.Range("B4:BC4").Value = $row4

# --- Row 5 ---
$row5 = New-Object 'object[,]' 1,54
$row5[0,0] = 45
$row5[0,1] = "Central Africa 1"
$row5[0,2] = "Congo, Dem. Rep."
$row5[0,3] = "COD"
$row5[0,4] = "Lower middle income"
$row5[0,5] = 13.5
$row5[0,6] = 0.4
$row5[0,7] = 42
$row5[0,8] = 0.06336919100000001
$row5[0,9] = 6
$row5[0,10] = 24.9109993
$row5[0,11] = 10.29500008
$row5[0,12] = 81.34899901999999
$row5[0,13] = 440
$row5[0,14] = 12.20400047
$row5[0,15] = 6.447000027
$row5[0,16] = 36.83218885
$row5[0,17] = 5.131664248
$row5[0,18] = 35917650630
$row5[0,19] = 9.470288097999999
$row5[0,20] = 0.03
$row5[0,21] = 3
$row5[0,22] = 101
$row5[0,23] = 46.57451961
$row5[0,24] = 10.6
$row5[0,25] = 3.707999945
$row5[0,26] = 41.976
$row5[0,27] = 4.553658283
$row5[0,28] = 0
$row5[0,29] = 0
$row5[0,30] = 0
$row5[0,31] = 0
$row5[0,32] = 0
$row5[0,33] = 0
$row5[0,34] = 0
$row5[0,35] = 0
$row5[0,36] = 0
$row5[0,37] = 0
$row5[0,38] = 0
$row5[0,39] = 0
$row5[0,40] = 0
$row5[0,41] = 0
$row5[0,42] = 0
$row5[0,43] = 0
$row5[0,44] = 0
$row5[0,45] = 1
$row5[0,46] = 1
$row5[0,47] = 0
$row5[0,48] = 0
$row5[0,49] = 0
$row5[0,50] = 1
$row5[0,51] = 0
$row5[0,52] = 0
$row5[0,53] = 0
$ws.Range("B5:BC5").Value = $row5

# --- Row 6 ---
$row6 = New-Object 'object[,]' 1,54
$row6[0,0] = 46
$row6[0,1] = "Central Africa 1"
$row6[0,2] = "Congo, Rep."
$row6[0,3] = "COG"
$row6[0,4] = "Low income"
$row6[0,5] = 51.86239243
$row6[0,6] = 19.23297764
$row6[0,7] = 69.46530914
$row6[0,8] = 0.635369293
$row6[0,9] = 10
$row6[0,10] = 5.625
$row6[0,11] = 5.31799984
$row6[0,12] = 38.22900009
$row6[0,13] = 2520
$row6[0,14] = 25
$row6[0,15] = 36.77099991
$row6[0,16] = 72.98675034
$row6[0,17] = 20.36515297
$row6[0,18] = 14177437982
$row6[0,19] = 6.779916158
$row6[0,20] = 0.29
$row6[0,21] = 7.11
$row6[0,22] = 56.5
$row6[0,23] = 53.17355571
$row6[0,24] = 7.4
$row6[0,25] = 9.998000145000001
$row6[0,26] = 64.95699999999999
$row6[0,27] = 3.146204713
$row6[0,28] = 0
$row6[0,29] = 0
$row6[0,30] = 0
$row6[0,31] = 0
$row6[0,32] = 0
$row6[0,33] = 0
$row6[0,34] = 0
$row6[0,35] = 0
$row6[0,36] = 0
$row6[0,37] = 0
$row6[0,38] = 0
$row6[0,39] = 0
$row6[0,40] = 0
$row6[0,41] = 0
$row6[0,42] = 0
$row6[0,43] = 0
$row6[0,44] = 0
$row6[0,45] = 1
$row6[0,46] = 1
$row6[0,47] = 0
$row6[0,48] = 0
$row6[0,49] = 0
$row6[0,50] = 1
$row6[0,51] = 0
$row6[0,52] = 0
$row6[0,53] = 0
$ws.Range("B6:BC6").Value = $row6

# --- Row 7 ---
$row7 = New-Object 'object[,]' 1,54
$row7[0,0] = 48
$row7[0,1] = "Central Africa 1"
$row7[0,2] = "Cote d'Ivoire"
$row7[0,3] = "CIV"
$row7[0,4] = "Lower middle income"
$row7[0,5] = 61.9
$row7[0,6] = 36.54494476
$row7[0,7] = 83.9561824
$row7[0,8] = 0.490206046
$row7[0,9] = 10
$row7[0,10] = 18.99900055
$row7[0,11] = 8.211000443
$row7[0,12] = 50.375
$row7[0,13] = 1460
$row7[0,14] = 5.986000061
$row7[0,15] = 43.63999939
$row7[0,16] = 36.6588334
$row7[0,17] = 1.240430526
$row7[0,18] = 35372603446
$row7[0,19] = 8.79407739
$row7[0,20] = 0.15
$row7[0,21] = 19.2742298
$row7[0,22] = 98.3
$row7[0,23] = 26.27654966
$row7[0,24] = 9.4
$row7[0,25] = 2.747999907
$row7[0,26] = 53.479
$row7[0,27] = 3.881971936
$row7[0,28] = 0
$row7[0,29] = 0
$row7[0,30] = 0
$row7[0,31] = 0
$row7[0,32] = 0
$row7[0,33] = 0
$row7[0,34] = 0
$row7[0,35] = 0
$row7[0,36] = 0
$row7[0,37] = 0
$row7[0,38] = 0
$row7[0,39] = 0
$row7[0,40] = 0
$row7[0,41] = 0
$row7[0,42] = 0
$row7[0,43] = 0
$row7[0,44] = 0
$row7[0,45] = 1
$row7[0,46] = 0
$row7[0,47] = 0
$row7[0,48] = 0
$row7[0,49] = 0
$row7[0,50] = 0
$row7[0,51] = 0
$row7[0,52] = 0
$row7[0,53] = 0
$ws.Range("B7:BC7").Value = $row7

# --- Row 8 ---
$row8 = New-Object 'object[,]' 1,54
$row8[0,0] = 61
$row8[0,1] = "Central Africa 1"
$row8[0,2] = "Equatorial Guinea"
$row8[0,3] = "GNQ"
$row8[0,4] = "Upper middle income"
$row8[0,5] = 67.05870819
$row8[0,6] = 51.01377073
$row8[0,7] = 91.37228394
$row8[0,8] = 4.733816529
$row8[0,9] = 6
$row8[0,10] = 16.01600075
$row8[0,11] = 5.980000019
$row8[0,12] = 54.96300125
$row8[0,13] = 13140
$row8[0,14] = 8.237000464999999
$row8[0,15] = 36.79999924
$row8[0,16] = 65.9632886
$row8[0,17] = 0.772319187
$row8[0,18] = 21736500713
$row8[0,19] = 0.415061836
$row8[0,20] = 0.49
$row8[0,21] = 18.86
$row8[0,22] = 97.09999999999999
$row8[0,23] = 46.69405811
$row8[0,24] = 24
$row8[0,25] = 5.494999886
$row8[0,26] = 39.756
$row8[0,27] = 4.514037661
$row8[0,28] = 0
$row8[0,29] = 0
$row8[0,30] = 0
$row8[0,31] = 0
$row8[0,32] = 0
$row8[0,33] = 0
$row8[0,34] = 0
$row8[0,35] = 0
$row8[0,36] = 0
$row8[0,37] = 0
$row8[0,38] = 0
$row8[0,39] = 0
$row8[0,40] = 0
$row8[0,41] = 0
$row8[0,42] = 0
$row8[0,43] = 0
$row8[0,44] = 0
$row8[0,45] = 1
$row8[0,46] = 1
$row8[0,47] = 0
$row8[0,48] = 0
$row8[0,49] = 0
$row8[0,50] = 0
$row8[0,51] = 0
$row8[0,52] = 0
$row8[0,53] = 0
$ws.Range("B8:BC8").Value = $row8

# --- Row 9 ---
$row9 = New-Object 'object[,]' 1,54
$row9[0,0] = 74
$row9[0,1] = "Central Africa 1"
$row9[0,2] = "Ghana"
$row9[0,3] = "GHA"
$row9[0,4] = "Lower middle income"
$row9[0,5] = 78.3
$row9[0,6] = 63
$row9[0,7] = 90.8
$row9[0,8] = 0.536533378
$row9[0,9] = 11
$row9[0,10] = 27.77700043
$row9[0,11] = 15.79800034
$row9[0,12] = 44.72000122
$row9[0,13] = 1590
$row9[0,14] = 14.10700035
$row9[0,15] = 41.17300034
$row9[0,16] = 39.52355867
$row9[0,17] = 8.604962551
$row9[0,18] = 39086625009
$row9[0,19] = 3.985865624
$row9[0,20] = 0.11
$row9[0,21] = 25.51773743
$row9[0,22] = 63.4
$row9[0,23] = 25.35736494
$row9[0,24] = 10.9
$row9[0,25] = 2.164000034
$row9[0,26] = 53.392
$row9[0,27] = 3.550496942
$row9[0,28] = 0
$row9[0,29] = 0
$row9[0,30] = 0
$row9[0,31] = 0
$row9[0,32] = 0
$row9[0,33] = 0
$row9[0,34] = 0
$row9[0,35] = 0
$row9[0,36] = 0
$row9[0,37] = 0
$row9[0,38] = 0
$row9[0,39] = 0
$row9[0,40] = 0
$row9[0,41] = 0
$row9[0,42] = 0
$row9[0,43] = 0
$row9[0,44] = 0
$row9[0,45] = 1
$row9[0,46] = 1
$row9[0,47] = 0
$row9[0,48] = 0
$row9[0,49] = 0
$row9[0,50] = 1
$row9[0,51] = 0
$row9[0,52] = 0
$row9[0,53] = 0
$ws.Range("B9:BC9").Value = $row9

# --- Row 10 ---
$row10 = New-Object 'object[,]' 1,54
$row10[0,0] = 101
$row10[0,1] = "Central Africa 1"
$row10[0,2] = "Kenya"
$row10[0,3] = "KEN"
$row10[0,4] = "Lower middle income"
$row10[0,5] = 36
$row10[0,6] = 12.6
$row10[0,7] = 68.40000000000001
$row10[0,8] = 0.310415314
$row10[0,9] = 12
$row10[0,10] = 53.45399857
$row10[0,11] = 25.95800018
$row10[0,12] = 37.63899994
$row10[0,13] = 1260
$row10[0,14] = 14.90200043
$row10[0,15] = 47.45899963
$row10[0,16] = 18.29698092
$row10[0,17] = 1.335986481
$row10[0,18] = 61448046802
$row10[0,19] = 5.357116778
$row10[0,20] = 0.3
$row10[0,21] = 16.5
$row10[0,22] = 53.5
$row10[0,23] = 16.52198912
$row10[0,24] = 19.1
$row10[0,25] = 11.66699982
$row10[0,26] = 25.197
$row10[0,27] = 4.304922007
$row10[0,28] = 0
$row10[0,29] = 0
$row10[0,30] = 0
$row10[0,31] = 0
$row10[0,32] = 0
$row10[0,33] = 0
$row10[0,34] = 0
$row10[0,35] = 0
$row10[0,36] = 0
$row10[0,37] = 0
$row10[0,38] = 0
$row10[0,39] = 0
$row10[0,40] = 0
$row10[0,41] = 0
$row10[0,42] = 0
$row10[0,43] = 0
$row10[0,44] = 0
$row10[0,45] = 0
$row10[0,46] = 0
$row10[0,47] = 0
$row10[0,48] = 0
$row10[0,49] = 0
$row10[0,50] = 0
$row10[0,51] = 0
$row10[0,52] = 0
$row10[0,53] = 0
$ws.Range("B10:BC10").Value = $row10

# --- Row 11 ---
$row11 = New-Object 'object[,]' 1,54
$row11[0,0] = 145
$row11[0,1] = "Central Africa 1"
$row11[0,2] = "Nigeria"
$row11[0,3] = "NGA"
$row11[0,4] = "Lower middle income"
$row11[0,5] = 56.37191391
$row11[0,6] = 31.67585522
$row11[0,7] = 84.2855835
$row11[0,8] = 0.545622113
$row11[0,9] = 9
$row11[0,10] = 12.18200016
$row11[0,11] = 9.477000237
$row11[0,12] = 36.77999878
$row11[0,13] = 2980
$row11[0,14] = 12.06200027
$row11[0,15] = 51.1590004
$row11[0,16] = 18.43512605
$row11[0,17] = 0.818201344
$row11[0,18] = 568499000000
$row11[0,19] = 6.309718596
$row11[0,20] = 0.21
$row11[0,21] = 21
$row11[0,22] = 111.6
$row11[0,23] = 41.60273017
$row11[0,24] = 6.7
$row11[0,25] = 4.559999943
$row11[0,26] = 46.942
$row11[0,27] = 4.48255153
$row11[0,28] = 0
$row11[0,29] = 0
$row11[0,30] = 0
$row11[0,31] = 0
$row11[0,32] = 0
$row11[0,33] = 0
$row11[0,34] = 0
$row11[0,35] = 0
$row11[0,36] = 0
$row11[0,37] = 0
$row11[0,38] = 0
$row11[0,39] = 0
$row11[0,40] = 0
$row11[0,41] = 0
$row11[0,42] = 0
$row11[0,43] = 0
$row11[0,44] = 0
$row11[0,45] = 1
$row11[0,46] = 1
$row11[0,47] = 0
$row11[0,48] = 0
$row11[0,49] = 0
$row11[0,50] = 1
$row11[0,51] = 0
$row11[0,52] = 0
$row11[0,53] = 0
$ws.Range("B11:BC11").Value = $row11

# --- Row 12 ---
$row12 = New-Object 'object[,]' 1,54
$row12[0,0] = 162
$row12[0,1] = "Central Africa 1"
$row12[0,2] = "Rwanda"
$row12[0,3] = "RWA"
$row12[0,4] = "Low income"
$row12[0,5] = 19.8
$row12[0,6] = 9.1
$row12[0,7] = 71.8
$row12[0,8] = 0.074016446
$row12[0,9] = 6
$row12[0,10] = 12.77900028
$row12[0,11] = 6.039000034
$row12[0,12] = 68.45300293
$row12[0,13] = 700
$row12[0,14] = 7.913000107
$row12[0,15] = 23.63400078
$row12[0,16] = 14.72211453
$row12[0,17] = 3.926137464
$row12[0,18] = 8016591928
$row12[0,19] = 7.62457575
$row12[0,20] = 0.14
$row12[0,21] = 10.6
$row12[0,22] = 43.4
$row12[0,23] = 51.86310993
$row12[0,24] = 63.8
$row12[0,25] = 1.169999957
$row12[0,26] = 27.841
$row12[0,27] = 6.054470586
$row12[0,28] = 0
$row12[0,29] = 0
$row12[0,30] = 0
$row12[0,31] = 0
$row12[0,32] = 0
$row12[0,33] = 0
$row12[0,34] = 0
$row12[0,35] = 0
$row12[0,36] = 0
$row12[0,37] = 0
$row12[0,38] = 0
$row12[0,39] = 0
$row12[0,40] = 0
$row12[0,41] = 0
$row12[0,42] = 0
$row12[0,43] = 0
$row12[0,44] = 0
$row12[0,45] = 1
$row12[0,46] = 1
$row12[0,47] = 0
$row12[0,48] = 0
$row12[0,49] = 0
$row12[0,50] = 0
$row12[0,51] = 0
$row12[0,52] = 0
$row12[0,53] = 0
$ws.Range("B12:BC12").Value = $row12

# --- Row 13 ---
$row13 = New-Object 'object[,]' 1,54
$row13[0,0] = 167
$row13[0,1] = "Central Africa 1"
$row13[0,2] = "Senegal"
$row13[0,3] = "SEN"
$row13[0,4] = "Low income"
$row13[0,5] = 61
$row13[0,6] = 32.7
$row13[0,7] = 85
$row13[0,8] = 0.608809117
$row13[0,9] = 11
$row13[0,10] = 16.31599998
$row13[0,11] = 5.46999979
$row13[0,12] = 53.67399979
$row13[0,13] = 1030
$row13[0,14] = 20.20499992
$row13[0,15] = 26.12100029
$row13[0,16] = 28.12071229
$row13[0,17] = 2.630373382
$row13[0,18] = 15304363138
$row13[0,19] = 4.075083326
$row13[0,20] = 0.01
$row13[0,21] = 17.7
$row13[0,22] = 52.3
$row13[0,23] = 35.99474681
$row13[0,24] = 43.3
$row13[0,25] = 6.356999874
$row13[0,26] = 43.393
$row13[0,27] = 3.697122779
$row13[0,28] = 0
$row13[0,29] = 0
$row13[0,30] = 0
$row13[0,31] = 0
$row13[0,32] = 0
$row13[0,33] = 0
$row13[0,34] = 0
$row13[0,35] = 0
$row13[0,36] = 0
$row13[0,37] = 0
$row13[0,38] = 0
$row13[0,39] = 0
$row13[0,40] = 0
$row13[0,41] = 0
$row13[0,42] = 0
$row13[0,43] = 0
$row13[0,44] = 0
$row13[0,45] = 1
$row13[0,46] = 1
$row13[0,47] = 0
$row13[0,48] = 0
$row13[0,49] = 0
$row13[0,50] = 0
$row13[0,51] = 0
$row13[0,52] = 0
$row13[0,53] = 0
$ws.Range("B13:BC13").Value = $row13

# --- Row 14 ---
$row14 = New-Object 'object[,]' 1,54
$row14[0,0] = 185
$row14[0,1] = "Central Africa 1"
$row14[0,2] = "Sudan"
$row14[0,3] = "SDN"
$row14[0,4] = "Lower middle income"
$row14[0,5] = 44.9
$row14[0,6] = 31.7
$row14[0,7] = 76.3
$row14[0,8] = 0.299732598
$row14[0,9] = 8
$row14[0,10] = 26.11599922
$row14[0,11] = 9.135000228999999
$row14[0,12] = 53.35599899
$row14[0,13] = 1830
$row14[0,14] = 18.83399963
$row14[0,15] = 27.80999947
$row14[0,16] = 8.149134518
$row14[0,17] = 1.523136574
$row14[0,18] = 82151588419
$row14[0,19] = 2.679411813
$row14[0,20] = 0.02
$row14[0,21] = 24.64
$row14[0,22] = 69
$row14[0,23] = 52.51987476
$row14[0,24] = 24.3
$row14[0,25] = 12.74800014
$row14[0,26] = 33.623
$row14[0,27] = 2.867143272
$row14[0,28] = 0
$row14[0,29] = 0
$row14[0,30] = 0
$row14[0,31] = 0
$row14[0,32] = 0
$row14[0,33] = 0
$row14[0,34] = 0
$row14[0,35] = 0
$row14[0,36] = 0
$row14[0,37] = 0
$row14[0,38] = 0
$row14[0,39] = 0
$row14[0,40] = 0
$row14[0,41] = 0
$row14[0,42] = 0
$row14[0,43] = 0
$row14[0,44] = 0
$row14[0,45] = 1
$row14[0,46] = 1
$row14[0,47] = 0
$row14[0,48] = 0
$row14[0,49] = 0
$row14[0,50] = 1
$row14[0,51] = 0
$row14[0,52] = 0
$row14[0,53] = 0
$ws.Range("B14:BC14").Value = $row14

# --- Row 15 ---
$row15 = New-Object 'object[,]' 1,54
$row15[0,0] = 203
$row15[0,1] = "Central Africa 1"
$row15[0,2] = "Uganda"
$row15[0,3] = "UGA"
$row15[0,4] = "Low income"
$row15[0,5] = 20.4
$row15[0,6] = 10.3
$row15[0,7] = 51.4
$row15[0,8] = 0.134656001
$row15[0,9] = 7
$row15[0,10] = 19.75
$row15[0,11] = 12.92700005
$row15[0,12] = 71.54699707
$row15[0,13] = 660
$row15[0,14] = 6.887000084
$row15[0,15] = 21.56500053
$row15[0,16] = 18.17025675
$row15[0,17] = 3.878679401
$row15[0,18] = 27291880327
$row15[0,19] = 5.106307324
$row15[0,20] = 0.38
$row15[0,21] = 16.9
$row15[0,22] = 60.1
$row15[0,23] = 61.38641667
$row15[0,24] = 35
$row15[0,25] = 1.907999992
$row15[0,26] = 15.766
$row15[0,27] = 5.459493493
$row15[0,28] = 0
$row15[0,29] = 0
$row15[0,30] = 0
$row15[0,31] = 0
$row15[0,32] = 0
$row15[0,33] = 0
$row15[0,34] = 0
$row15[0,35] = 0
$row15[0,36] = 0
$row15[0,37] = 0
$row15[0,38] = 0
$row15[0,39] = 0
$row15[0,40] = 0
$row15[0,41] = 0
$row15[0,42] = 0
$row15[0,43] = 0
$row15[0,44] = 0
$row15[0,45] = 0
$row15[0,46] = 1
$row15[0,47] = 0
$row15[0,48] = 0
$row15[0,49] = 0
$row15[0,50] = 1
$row15[0,51] = 0
$row15[0,52] = 0
$row15[0,53] = 0
$ws.Range("B15:BC15").Value = $row15

